# Apply updated "dSF" (column F) values for a set of rows on Sheet1.
# This mirrors a repull/push of the underlying data plus a recalculation
# of the mean (dSF) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$updates = @{
    10 = 2
    13 = -2
    20 = 4
    22 = -3
    29 = 0
    30 = -2
    34 = 1
    35 = 5
    36 = 2
    37 = -4
    43 = 0
    45 = -2
    56 = -9
    61 = -2
    67 = -1
    70 = -5
    74 = 0
    76 = 4
    85 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
